$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Test scenario / test case headers
$ws.Range("A2").Value = "TestScenario_1"
$ws.Range("B2").Value = "Create_ Account_01"

# Step data for rows 2-8 in columns F (step no), G (user action), H (expected result)
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "Url"
$ws.Range("H2").Value = "Url should be launched"

$ws.Range("F3").Value = "2"
$ws.Range("G3").Value = "Enter UserName"
$ws.Range("H3").Value = "User Name should be entered"

$ws.Range("F4").Value = "3"
$ws.Range("G4").Value = "Enter Password"
$ws.Range("H4").Value = "Password should be entered"

$ws.Range("F5").Value = "4"
$ws.Range("G5").Value = "Click SignIn"
$ws.Range("H5").Value = "Sign in should be clicked"

$ws.Range("F6").Value = "5"
$ws.Range("G6").Value = "Click on Case"
$ws.Range("H6").Value = "Accounts page will get opens"

$ws.Range("F7").Value = "6"
$ws.Range("G7").Value = "Close"
$ws.Range("H7").Value = "PopUp Close"

$ws.Range("F8").Value = "7"
$ws.Range("G8").Value = "Verify"
$ws.Range("H8").Value = "Dialogue box should get close"

# Resize the table to include the new rows
$wb.Worksheets.Item(1).ListObjects.Item("Table1").Resize($ws.Range("A1:J8"))

# Adjust column widths to match target
$ws.Columns.Item(2).ColumnWidth = 20.410625
$ws.Columns.Item(7).ColumnWidth = 17.410625
$ws.Columns.Item(8).ColumnWidth = 29.410625
